$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update prepends a new record as row 7 (the most recent date),
# pushing the previously-existing rows 7-16 down to rows 8-17.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with this week's record.
$ws.Cells.Item(7, 1).Value2  = 8
$ws.Cells.Item(7, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(7, 3).Value2  = "Coquimbo"
$ws.Cells.Item(7, 4).Value2  = 44495
$ws.Cells.Item(7, 5).Value2  = 4
$ws.Cells.Item(7, 6).Value2  = "Fruta"
$ws.Cells.Item(7, 7).Value2  = 100101
$ws.Cells.Item(7, 8).Value2  = "Berries"
$ws.Cells.Item(7, 9).Value2  = 100101001
$ws.Cells.Item(7, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(7, 11).Value2 = "Sin especificar"
$ws.Cells.Item(7, 12).Value2 = "Primera"
$ws.Cells.Item(7, 13).Value2 = 300
$ws.Cells.Item(7, 14).Value2 = 11000
$ws.Cells.Item(7, 15).Value2 = 12000
$ws.Cells.Item(7, 16).Value2 = 11500
$ws.Cells.Item(7, 17).Value2 = "$/bandeja 2 kilos"
$ws.Cells.Item(7, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(7, 19).Value2 = 5750
$ws.Cells.Item(7, 20).Value2 = 2
